$d = $word.ActiveDocument
$d.Content.Find.Execute("722×6=", $true, $false, $false, $false, $false, $true, 1, $false, "554×2=", 2)
$d.Content.Find.Execute("546×6=", $true, $false, $false, $false, $false, $true, 1, $false, "467×3=", 2)
$d.Content.Find.Execute("720×9=", $true, $false, $false, $false, $false, $true, 1, $false, "683×3=", 2)
$d.Content.Find.Execute("912×5=", $true, $false, $false, $false, $false, $true, 1, $false, "716×4=", 2)
$d.Content.Find.Execute("617×5=", $true, $false, $false, $false, $false, $true, 1, $false, "835×6=", 2)
$d.Content.Find.Execute("742×7=", $true, $false, $false, $false, $false, $true, 1, $false, "764×2=", 2)
$d.Content.Find.Execute("897×6=", $true, $false, $false, $false, $false, $true, 1, $false, "578×7=", 2)
$d.Content.Find.Execute("219×3=", $true, $false, $false, $false, $false, $true, 1, $false, "890×4=", 2)
$d.Content.Find.Execute("876×7=", $true, $false, $false, $false, $false, $true, 1, $false, "775×3=", 2)
$d.Content.Find.Execute("140×5=", $true, $false, $false, $false, $false, $true, 1, $false, "872×6=", 2)
$d.Content.Find.Execute("814×3=", $true, $false, $false, $false, $false, $true, 1, $false, "583×6=", 2)
$d.Content.Find.Execute("868×4=", $true, $false, $false, $false, $false, $true, 1, $false, "141×3=", 2)
$d.Content.Find.Execute("983×6=", $true, $false, $false, $false, $false, $true, 1, $false, "337×5=", 2)
$d.Content.Find.Execute("562×5=", $true, $false, $false, $false, $false, $true, 1, $false, "473×5=", 2)
$d.Content.Find.Execute("363×2=", $true, $false, $false, $false, $false, $true, 1, $false, "452×9=", 2)
$d.Content.Find.Execute("524×5=", $true, $false, $false, $false, $false, $true, 1, $false, "478×8=", 2)
$d.Content.Find.Execute("211×9=", $true, $false, $false, $false, $false, $true, 1, $false, "473×8=", 2)
$d.Content.Find.Execute("492×4=", $true, $false, $false, $false, $false, $true, 1, $false, "113×6=", 2)
$d.Content.Find.Execute("247×5=", $true, $false, $false, $false, $false, $true, 1, $false, "475×9=", 2)
$d.Content.Find.Execute("460×4=", $true, $false, $false, $false, $false, $true, 1, $false, "950×3=", 2)
$d.Content.Find.Execute("556×2=", $true, $false, $false, $false, $false, $true, 1, $false, "593×3=", 2)
$d.Content.Find.Execute("429×5=", $true, $false, $false, $false, $false, $true, 1, $false, "559×7=", 2)
$d.Content.Find.Execute("593×2=", $true, $false, $false, $false, $false, $true, 1, $false, "131×3=", 2)
$d.Content.Find.Execute("929×7=", $true, $false, $false, $false, $false, $true, 1, $false, "916×4=", 2)
$d.Content.Find.Execute("457×4=", $true, $false, $false, $false, $false, $true, 1, $false, "506×8=", 2)
